$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Unprotect()

# Update the confidentiality/disclaimer note date from 2021-05-25 to 2021-05-26
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-26 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) values for rows 2-6
$ws.Range("D2").Value = 0.2543423447170707
$ws.Range("E2").Value = 0.0003463003578436652

$ws.Range("D3").Value = 0.2532819483341827
$ws.Range("E3").Value = 0.002942750133761329

$ws.Range("D4").Value = 0.2433113877418613
$ws.Range("E4").Value = 0.00319241559446648

$ws.Range("D5").Value = 0.2490643192068853
$ws.Range("E5").Value = 0.006518404907975395

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0.003233677481588204

$ws.Protect()
